# feat: add 2022-Q3 data
#
# 1. "总计" (summary) sheet: push the existing quarter rows down by one and
#    insert a new 2022-Q3 summary row at the top of the data block.
# 2. Insert a brand-new worksheet named "2022-Q3" right after "总计" holding
#    the per-fund holding detail for that quarter (same layout as the other
#    quarterly sheets).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Part 1: "总计" summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Existing rows 2..6 (2022-Q2 .. 2020-Q4) move down to rows 3..7 (value-only
# shift -- the A-column index sequence 0,1,2,3,4 stays put and simply grows
# one more entry at the bottom).
$existing = @()
for ($r = 2; $r -le 6; $r++) {
    $existing += , @($summary.Cells.Item($r, 2).Value2, $summary.Cells.Item($r, 3).Value2, $summary.Cells.Item($r, 4).Value2)
}

for ($i = 0; $i -lt $existing.Length; $i++) {
    $targetRow = $i + 3
    $vals = $existing[$i]
    $summary.Cells.Item($targetRow, 2).Value = $vals[0]
    $summary.Cells.Item($targetRow, 3).Value = $vals[1]
    $summary.Cells.Item($targetRow, 4).Value = $vals[2]
}

# New row 7 needs the same index-column styling as the other data rows --
# clone it from row 6's A cell so it picks up the identical style record.
$summary.Cells.Item(6, 1).Copy($summary.Cells.Item(7, 1))
$summary.Cells.Item(7, 1).Value = 5

# New 2022-Q3 summary row at row 2.
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 12
$summary.Cells.Item(2, 4).Value = 0.78

# ---------------------------------------------------------------------------
# Part 2: brand-new "2022-Q3" worksheet (per-fund detail)
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

# Match the page-margin defaults used by the other quarterly sheets.
$q3.PageSetup.LeftMargin = 54
$q3.PageSetup.RightMargin = 54
$q3.PageSetup.TopMargin = 72
$q3.PageSetup.BottomMargin = 72
$q3.PageSetup.HeaderMargin = 36
$q3.PageSetup.FooterMargin = 36

# Clone header-row formatting (bold / centered / bordered) from the summary
# sheet's own header cells, then overwrite with the Q3 sheet's headers.
$summary.Range("B1:D1").Copy($q3.Range("B1:H1"))

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q3.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$rows = @(
    @("010296", "万家互联互通中国优势量化策略混合A", "4.22", "86.53", "5.46", "0.2304", 1),
    @("001305", "九泰天富改革新动力混合A", "1.86", "94.71", "7.98", "0.1484", 5),
    @("004823", "上投摩根安裕回报混合A", "5.35", "25.71", "1.62", "0.0867", 6),
    @("004824", "上投摩根安裕回报混合C", "4.91", "25.71", "1.62", "0.0795", 6),
    @("001844", "九泰久益灵活配置混合C", "0.98", "93.32", "7.71", "0.0756", 7),
    @("206013", "鹏华宏观灵活配置混合", "0.95", "72.10", "4.65", "0.0442", 7),
    @("001782", "九泰久益灵活配置混合A", "0.53", "93.32", "7.71", "0.0409", 7),
    @("014938", "同泰产业升级混合A", "1.01", "61.58", "2.77", "0.0280", 8),
    @("010297", "万家互联互通中国优势量化策略混合C", "0.46", "86.53", "5.46", "0.0251", 1),
    @("009912", "九泰天富改革新动力混合C", "0.17", "94.71", "7.98", "0.0136", 5),
    @("562530", "华夏中证智选1000价值稳健策略ETF", "0.54", "94.32", "0.96", "0.0052", 2),
    @("014939", "同泰产业升级混合C", "0.00", "61.58", "2.77", 0, 8)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    # Index column (A) -- bold/centered/bordered, cloned from the summary
    # sheet's index column so the style record matches exactly.
    $summary.Cells.Item(2, 1).Copy($q3.Cells.Item($r, 1))
    $q3.Cells.Item($r, 1).Value = $i

    $q3.Cells.Item($r, 2).NumberFormat = "@"
    $q3.Cells.Item($r, 2).Value = $row[0]

    $q3.Cells.Item($r, 3).NumberFormat = "@"
    $q3.Cells.Item($r, 3).Value = $row[1]

    $q3.Cells.Item($r, 4).NumberFormat = "@"
    $q3.Cells.Item($r, 4).Value = $row[2]

    $q3.Cells.Item($r, 5).NumberFormat = "@"
    $q3.Cells.Item($r, 5).Value = $row[3]

    $q3.Cells.Item($r, 6).NumberFormat = "@"
    $q3.Cells.Item($r, 6).Value = $row[4]

    if ($i -eq ($rows.Length - 1)) {
        # last row's market-value column is a genuine 0, stored as a number
        $q3.Cells.Item($r, 7).Value = $row[5]
    } else {
        $q3.Cells.Item($r, 7).NumberFormat = "@"
        $q3.Cells.Item($r, 7).Value = $row[5]
    }

    $q3.Cells.Item($r, 8).Value = $row[6]
}
